$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity Log - Part 2")

# Row 19
$ws.Cells.Item(19, 2).Value = 624
$ws.Cells.Item(19, 3).Value = 43937
$ws.Cells.Item(19, 4).Value = 0.66666666666666663
$ws.Cells.Item(19, 5).Value = 0.73611111111111116
$ws.Cells.Item(19, 7).Value = "Deon with Documenting Design Entities for all our designs"

# Row 20
$ws.Cells.Item(20, 2).Value = 624
$ws.Cells.Item(20, 3).Value = 43937
$ws.Cells.Item(20, 4).Value = 0.73611111111111116
$ws.Cells.Item(20, 5).Value = 0.78333333333333333
$ws.Cells.Item(20, 7).Value = "Adding Comments to Source code"

# Row 21
$ws.Cells.Item(21, 2).Value = 624
$ws.Cells.Item(21, 3).Value = 43937
$ws.Cells.Item(21, 4).Value = 0.83333333333333337
$ws.Cells.Item(21, 5).Value = 0.89930555555555547
$ws.Cells.Item(21, 7).Value = "Worked on improving report"

# Row 22
$ws.Cells.Item(22, 2).Value = 624
$ws.Cells.Item(22, 3).Value = 43937
$ws.Cells.Item(22, 4).Value = 0.89930555555555547
$ws.Cells.Item(22, 5).Value = 0.95138888888888884
$ws.Cells.Item(22, 7).Value = "Done with discussion on the observation and results of our functional simulation"

# Row 23
$ws.Cells.Item(23, 2).Value = 624
$ws.Cells.Item(23, 3).Value = 43938
$ws.Cells.Item(23, 4).Value = 0.67708333333333337
$ws.Cells.Item(23, 5).Value = 0.70277777777777783
$ws.Cells.Item(23, 7).Value = "Did a quick proofread and setup the folder for our documentation"

# Update the view: scroll position and selection on this sheet (topLeftCell A7, selection G23)
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("G23").Select()

# Update workbook window position (xWindow/yWindow)
$win.Left = 1725
$win.Top = 1800
